# Update generated statistics (view counts / min prices) on the
# "展览" and "全部类型" sheets to match the newly scraped data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 5636
    $ws.Range("F3").Value = 379
    $ws.Range("F5").Value = 311
    $ws.Range("F6").Value = 835
    $ws.Range("G6").Value = 75
    $ws.Range("F7").Value = 56

    if ($name -eq "展览") {
        $ws.Range("F8").Value = 374
    } else {
        $ws.Range("F9").Value = 374
    }
}
